$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $idx = -1
            for ($i = 0; $i -lt $parts.Count; $i++) {
                if ($parts[$i] -ceq "System") {
                    $idx = $i
                }
            }
            if ($idx -ge 0 -and $idx -ne ($parts.Count - 1)) {
                $newParts = @()
                for ($i = 0; $i -lt $parts.Count; $i++) {
                    if ($i -ne $idx) {
                        $newParts += $parts[$i]
                    }
                }
                $newParts += "System"
                $newVal = $newParts -join ", "
                $cell.Value = $newVal
            }
        }
    }
}
